# 1st changes of mifos to finflux
#
# Insert a new (blank) column before column N ("Late") on the
# "Repayment schedule" sheet, shifting the old N/O/P columns to O/P/Q,
# and make "Repayment schedule" the active/selected sheet (instead of
# "ModifyLoan").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at position N (14) - pushes Late/heading/Outstanding
# columns one slot to the right.
$ws.Columns.Item(14).Insert() | Out-Null

# Give the freshly inserted column the same width the old column N/"Late"
# column used to have.
$ws.Columns.Item(14).ColumnWidth = 9.86

# Make "Repayment schedule" the active sheet/tab and update its selection,
# which also clears the previously selected cell on "ModifyLoan".
$ws.Activate() | Out-Null
$ws.Range("S6").Select() | Out-Null
